# Update "Generate Report for Handback" timestamps.
$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 27486164-... (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-26 00:45:01"

# --- zh-cn sheet: row for 27486164-... (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-26 00:44:54"
$wsZhCn.Range("K3").Value = "2016-08-26 00:45:30"

# --- de-de sheet: row for 27486164-... (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-26 00:45:37"
